$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.426.31'
$ws.Range("E2").Value = '  -1.73%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.620.76'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '533.18'
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.27'
$ws.Range("E6").Value = '  +0.66%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.566'
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.93'
$ws.Range("E9").Value = '  +6.85%  '
$ws.Range("E10").Value = '  -1.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.333'
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.089.36'
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '58.371.51'
$ws.Range("E14").Value = '  -1.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.69'
$ws.Range("E15").Value = '  -0.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.631.63'
$ws.Range("E16").Value = '  +0.93%  '
$ws.Range("E17").Value = '  -1.20%  '
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '333.93'
$ws.Range("E19").Value = '  -2.17%  '
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("E21").Value = '  -2.25%  '
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.29'
$ws.Range("E23").Value = '  -1.83%  '
$ws.Range("E24").Value = '  +1.63%  '
$ws.Range("B25").Value = 'Binance-PegBSC-USD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.162'
$ws.Range("E26").Value = '  -1.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.08'
$ws.Range("E27").Value = '  -2.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0734'
$ws.Range("E28").Value = '  -1.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  -1.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.85'
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.75'
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '150.21'
$ws.Range("E33").Value = '  +0.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.88'
$ws.Range("E34").Value = '  -2.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.849'
$ws.Range("E35").Value = '  +1.52%  '
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("E37").Value = '  -3.84%  '
$ws.Range("E38").Value = '  -2.04%  '
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '280.07'
$ws.Range("E40").Value = '  +3.09%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.594'
$ws.Range("E42").Value = '  -0.91%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.69'
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.95'
$ws.Range("E44").Value = '  +2.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0527'
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0935'
$ws.Range("E46").Value = '  -1.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0223'
$ws.Range("E47").Value = '  +0.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.937.16'
$ws.Range("E48").Value = '  -0.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.43'
$ws.Range("E49").Value = '  -1.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.80'
$ws.Range("E50").Value = '  -4.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '112.71'
$ws.Range("E51").Value = '  +1.43%  '
